{"js": "// Append a new paragraph at the end of the document body (after the\n// \"${foo}\" paragraph, before the sectPr) that demonstrates resolving an\n// expression that calls a method on a Map: ${map.get(\"xyz\")}.\nconst body = context.document.body;\n\nconst newPara = body.insertParagraph(\n  'This is a more advanced object: ${map.get(\"xyz\")}.',\n  Word.InsertLocation.end\n);\n\n// Match the spacing used by the sibling \"${foo}\" paragraph\n// (w:spacing w:before=\"0\" w:after=\"140\" -> points are twentieths of a point).\nnewPara.spaceBefore = 0;\nnewPara.spaceAfter = 7;\n\nawait context.sync();\n", "ps1": "# Append a new paragraph at the end of the document (after the \"${foo}\"\n# paragraph) that shows resolving an expression calling a method on a Map:\n# ${map.get(\"xyz\")}.\n$d = $word.ActiveDocument\n\n$lastIndex = $d.Paragraphs.Count\n$lastPara = $d.Paragraphs.Item($lastIndex)\n\n# Splits a new paragraph off right after the current last paragraph; the\n# new paragraph inherits the style (\"Brdtekst\") and run formatting\n# (w:lang=\"de-DE\") of the paragraph it was split from.\n$lastPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($d.Paragraphs.Count)\n$newPara.Range.Text = 'This is a more advanced object: ${map.get(\"xyz\")}.'\n\n# Match the spacing used by the sibling \"${foo}\" paragraph\n# (w:spacing w:before=\"0\" w:after=\"140\" -- points are twentieths of a point).\n$newPara.SpaceBefore = 0\n$newPara.SpaceAfter = 7\n"}
